$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2: dimension/measure identifiers curated
$ws.Range("A2").Value = "iaest-measure:horas-trabajadas"
$ws.Range("C2").Value = "sdmx-dimension:refArea"
$ws.Range("D2").Value = "iaest-measure:situacion-profesional"

# Update row 3: type labels (dim/medida) swapped per column
$ws.Range("A3").Value = "medida"
$ws.Range("B3").Value = "medida"
$ws.Range("C3").Value = "dim"
$ws.Range("D3").Value = "medida"

# Update row 4: data type / concept scheme values
$ws.Range("A4").Value = "xsd:int"
$ws.Range("B4").Value = "xsd:int"
$ws.Range("C4").Value = "URI-Comunidad"
$ws.Range("D4").Value = "xsd:int"

# Row 5 (mapping file references) is no longer needed
$ws.Rows.Item(5).Delete()
